# Fix the "MasterCart" -> "MasterCard" typo.
# (commit: "arreglado un typo en la presentación que la hacía LITERALLY
# UNPLAYABLE")
#
# Before (3 runs):  " Visa, "            "MasterCart"  " y American Express"
# After  (3 runs):  " Visa"              ", MasterCard "  "y American Express"
#
# We find the shape containing the typo by scanning the deck (instead of
# hard-coding a slide/shape index) and rewrite only the affected substring
# via TextRange.Characters, so every other run in the paragraph/shape is
# left completely untouched.

$p = $ppt.ActivePresentation

$fullSpan = " Visa, MasterCart y American Express"

foreach ($s in $p.Slides) {
    foreach ($sh in $s.Shapes) {
        if (-not $sh.HasTextFrame) { continue }
        if (-not $sh.TextFrame.HasText) { continue }

        $tr = $sh.TextFrame.TextRange
        $spanPos = $tr.Text.IndexOf($fullSpan)
        if ($spanPos -lt 0) { continue }

        $base1 = $spanPos + 1   # 1-based start of $fullSpan for Characters()

        # Segment layout inside $fullSpan (lengths in characters):
        #   " Visa"               ->  5  (untouched)
        #   ", MasterCart "       -> 13  (text corrected to ", MasterCard ")
        #   "y American Express"  -> 18  (untouched)
        $typo = $tr.Characters($base1 + 5, 13)
        $typo.Text = ", MasterCard "

        break
    }
}
